# Updates cryptos list prices / volume(1h) percentages (and fixes the
# WrappedBTC / WrappedEther row order) to match the latest scrape.
# For cells whose new text looks like a plain number (e.g. "0.529"),
# a leading apostrophe is used so Excel keeps storing them as text
# (matching the original inlineStr string cells) instead of silently
# converting them to numeric cells; the style is then reset back to
# "Normal" so no stray number-format/quote-prefix styling is left on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.776.73"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "3.143.30"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'587.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'145.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.140.17"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.71%  "
$ws.Range("D11").Value = "'5.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'37.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").Value = "3.663.58"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.144.48"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.585.77"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'463.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'14.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "'0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").Value = "'13.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "'81.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'9.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.75%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'6.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'26.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "'0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "0.0₃0854"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").Value = "'6.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "'50.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").Value = "'440.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "'8.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "2.916.21"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'0.279"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "'36.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.95%  "
$ws.Range("D48").Value = "'125.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "'24.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
